$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value for existing row 19 (X-n251-k28) -> Maximal Path Length column
$ws.Range("D19").Value = "36917,9…"

# New row 21: X-n289-k60
$ws.Range("A21").Value = "X-n289-k60"

# Row 22: X-n502-k39 (first result) - filled in the order the author typed them
$ws.Range("B22").Value = "69226(BKS)"
$ws.Range("B20").Value = "18839(BKS)"
$ws.Range("D22").Value = "67560,2…"
$ws.Range("E22").Value = "262s mit Startheuristik"
$ws.Range("A22").Value = "X-n502-k39"
$ws.Range("F22").Value = "68165,5… (evt. Mit SPPRC Pfaden)"

# Remaining numeric cells
$ws.Range("B21").Value = 95151
$ws.Range("C21").Value = 269
$ws.Range("C22").Value = 15

# Column E needs to be widened to fit the new long text (bestFit-style width)
$ws.Columns("E:E").ColumnWidth = 19.166666666666668

# Update the active selection to reflect where the author left off
$ws.Range("E23").Select() | Out-Null
